# Update automàtic: dades i banners [2026-02-18 11:50]
# Updates DATA_EXTRACCIO (column E) timestamps for each station row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 11:48:25"
$ws.Range("E3").Value = "2026-02-18 11:48:28"
$ws.Range("E4").Value = "2026-02-18 11:48:30"
$ws.Range("E5").Value = "2026-02-18 11:48:32"
$ws.Range("E6").Value = "2026-02-18 11:48:35"
$ws.Range("E7").Value = "2026-02-18 11:48:37"
$ws.Range("E8").Value = "2026-02-18 11:48:40"
$ws.Range("E9").Value = "2026-02-18 11:48:42"
$ws.Range("E10").Value = "2026-02-18 11:48:45"
$ws.Range("E11").Value = "2026-02-18 11:48:46"
$ws.Range("E12").Value = "2026-02-18 11:48:47"
$ws.Range("E13").Value = "2026-02-18 11:48:48"
$ws.Range("E14").Value = "2026-02-18 11:48:49"
$ws.Range("E15").Value = "2026-02-18 11:48:51"
$ws.Range("E16").Value = "2026-02-18 11:48:54"
$ws.Range("E17").Value = "2026-02-18 11:48:56"
$ws.Range("E18").Value = "2026-02-18 11:48:58"
$ws.Range("E19").Value = "2026-02-18 11:49:01"
$ws.Range("E20").Value = "2026-02-18 11:49:03"
$ws.Range("E21").Value = "2026-02-18 11:49:06"
$ws.Range("E22").Value = "2026-02-18 11:49:08"
$ws.Range("E23").Value = "2026-02-18 11:49:11"
$ws.Range("E24").Value = "2026-02-18 11:49:13"
$ws.Range("E25").Value = "2026-02-18 11:49:16"
$ws.Range("E26").Value = "2026-02-18 11:49:18"
$ws.Range("E27").Value = "2026-02-18 11:49:20"
$ws.Range("E28").Value = "2026-02-18 11:49:22"
$ws.Range("E29").Value = "2026-02-18 11:49:25"
$ws.Range("E30").Value = "2026-02-18 11:49:27"
$ws.Range("E31").Value = "2026-02-18 11:49:30"
$ws.Range("E32").Value = "2026-02-18 11:49:32"
$ws.Range("E33").Value = "2026-02-18 11:49:34"
$ws.Range("E34").Value = "2026-02-18 11:49:37"
$ws.Range("E35").Value = "2026-02-18 11:49:39"
$ws.Range("E36").Value = "2026-02-18 11:49:42"
$ws.Range("E37").Value = "2026-02-18 11:49:44"
$ws.Range("E38").Value = "2026-02-18 11:49:46"
$ws.Range("E39").Value = "2026-02-18 11:49:49"
$ws.Range("E40").Value = "2026-02-18 11:49:51"
$ws.Range("E41").Value = "2026-02-18 11:49:54"
$ws.Range("E42").Value = "2026-02-18 11:49:56"
$ws.Range("E43").Value = "2026-02-18 11:49:59"
$ws.Range("E44").Value = "2026-02-18 11:50:01"
$ws.Range("E45").Value = "2026-02-18 11:50:03"
$ws.Range("E46").Value = "2026-02-18 11:50:06"
